$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$setting = "100 english words wiht letter [A]"

$rows = @(
    @{A="28.02.2025_12.26.53"; B="en - de"; C="account";   D="Konto";        E="Konto";        F="ok";  G=$setting; H=""},
    @{A="28.02.2025_12.27.04"; B="en - de"; C="advertise"; D="werben";       E="inserat";      F="nok"; G=$setting; H=""},
    @{A="28.02.2025_12.27.15"; B="en - de"; C="acquire";   D="erwerben";     E="aquirieren";   F="nok"; G=$setting; H=""},
    @{A="28.02.2025_12.27.26"; B="en - de"; C="across";    D="über";         E="überqueren";   F="nok"; G=$setting; H=""},
    @{A="28.02.2025_12.27.38"; B="en - de"; C="alter";     D="ändern";       E="wechseln";     F="nok"; G=$setting; H=""},
    @{A="28.02.2025_12.27.45"; B="en - de"; C="alone";     D="allein";       E="allein";       F="ok";  G=$setting; H=""},
    @{A="28.02.2025_12.27.52"; B="en - de"; C="alone";     D="allein";       E="allein";       F="ok";  G=$setting; H=""},
    @{A="28.02.2025_12.27.58"; B="en - de"; C="all";       D="alle";         E="alle";         F="ok";  G=$setting; H=""},
    @{A="28.02.2025_12.28.07"; B="en - de"; C="allocate";  D="Zuteilen";     E="verbleiben";   F="nok"; G=$setting; H=""},
    @{A="28.02.2025_12.28.18"; B="en - de"; C="actual";    D="tatsächlich";  E="eigentlich";   F="nok"; G=$setting; H=""}
)

$startRow = 23
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}
